$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.297.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").Value = "'1.854.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.81%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'324.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4551"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.82%  "
$ws.Range("D8").Value = "'0.3878"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.76%  "
$ws.Range("E9").Value = "  -8.69%  "
$ws.Range("D10").Value = "'0.07943"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.09%  "
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("E12").Value = "  -4.22%  "
$ws.Range("D13").Value = "'1.851.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.99%  "
$ws.Range("D14").Value = "'5.898"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("E15").Value = "  -5.13%  "
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'85.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.10%  "
$ws.Range("D18").Value = "'0.06586"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'0.00001028"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("E20").Value = "  -5.96%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "'5.499"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.61%  "
$ws.Range("D23").Value = "'27.308.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").Value = "'10.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.63%  "
$ws.Range("D25").Value = "'2.290"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'2.086.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.93%  "
$ws.Range("D27").Value = "'153.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "'19.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").Value = "'2.060"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("D30").Value = "'5.478"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'121.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").Value = "'0.09342"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").Value = "'0.9367"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.47%  "
$ws.Range("D34").Value = "'1.457"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").Value = "'3.589"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").Value = "'5.262"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.54%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06020"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02227"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("D39").Value = "'1.222"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").Value = "'8.059"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.42%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'0.5918"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").Value = "'0.1884"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").Value = "'10.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.38%  "
$ws.Range("D45").Value = "'1.282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "'0.5613"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").Value = "'12.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.23%  "
$ws.Range("D48").Value = "'3.371"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("D49").Value = "'1.919"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.89%  "
$ws.Range("D50").Value = "'0.06734"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'108.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
